$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.10312032699585
$ws.Range("B1").Value = 1.902527809143066
$ws.Range("C1").Value = 9.266632080078125
$ws.Range("D1").Value = 2.402068614959717
$ws.Range("E1").Value = 1.289701342582703
